$d = $word.ActiveDocument
$find = $d.Content.Find
$find.Execute("2025-10-26 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-27 Monday", 2) | Out-Null
$find.Execute("82-77=", $true, $false, $false, $false, $false, $true, 1, $false, "90-87=", 2) | Out-Null
$find.Execute("33+35=", $true, $false, $false, $false, $false, $true, 1, $false, "35+51=", 2) | Out-Null
$find.Execute("98-45=", $true, $false, $false, $false, $false, $true, 1, $false, "87-21=", 2) | Out-Null
$find.Execute("35-0=", $true, $false, $false, $false, $false, $true, 1, $false, "87-31=", 2) | Out-Null
$find.Execute("53+6=", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=", 2) | Out-Null
$find.Execute("11+59=", $true, $false, $false, $false, $false, $true, 1, $false, "65+6=", 2) | Out-Null
$find.Execute("21+21=", $true, $false, $false, $false, $false, $true, 1, $false, "51-0=", 2) | Out-Null
$find.Execute("58-13=", $true, $false, $false, $false, $false, $true, 1, $false, "78-54=", 2) | Out-Null
$find.Execute("31+12=", $true, $false, $false, $false, $false, $true, 1, $false, "92-0=", 2) | Out-Null
$find.Execute("32+40=", $true, $false, $false, $false, $false, $true, 1, $false, "81+10=", 2) | Out-Null
$find.Execute("5+1=", $true, $false, $false, $false, $false, $true, 1, $false, "8+56=", 2) | Out-Null
$find.Execute("39+2=", $true, $false, $false, $false, $false, $true, 1, $false, "57-56=", 2) | Out-Null
$find.Execute("62-28=", $true, $false, $false, $false, $false, $true, 1, $false, "58-29=", 2) | Out-Null
$find.Execute("19+34=", $true, $false, $false, $false, $false, $true, 1, $false, "97-16=", 2) | Out-Null
$find.Execute("78-31=", $true, $false, $false, $false, $false, $true, 1, $false, "86+12=", 2) | Out-Null
$find.Execute("16+80=", $true, $false, $false, $false, $false, $true, 1, $false, "28+49=", 2) | Out-Null
$find.Execute("60-54=", $true, $false, $false, $false, $false, $true, 1, $false, "13+70=", 2) | Out-Null
$find.Execute("8+82=", $true, $false, $false, $false, $false, $true, 1, $false, "39+31=", 2) | Out-Null
$find.Execute("97+1=", $true, $false, $false, $false, $false, $true, 1, $false, "90-27=", 2) | Out-Null
$find.Execute("54-24=", $true, $false, $false, $false, $false, $true, 1, $false, "70+18=", 2) | Out-Null
$find.Execute("8+66=", $true, $false, $false, $false, $false, $true, 1, $false, "73-71=", 2) | Out-Null
$find.Execute("71+10=", $true, $false, $false, $false, $false, $true, 1, $false, "72-45=", 2) | Out-Null
$find.Execute("9+60=", $true, $false, $false, $false, $false, $true, 1, $false, "65+33=", 2) | Out-Null
$find.Execute("4+12=", $true, $false, $false, $false, $false, $true, 1, $false, "68-46=", 2) | Out-Null
$find.Execute("44-23=", $true, $false, $false, $false, $false, $true, 1, $false, "69-10=", 2) | Out-Null
$find.Execute("81+6=", $true, $false, $false, $false, $false, $true, 1, $false, "18+0=", 2) | Out-Null
$find.Execute("32-2=", $true, $false, $false, $false, $false, $true, 1, $false, "21+1=", 2) | Out-Null
$find.Execute("70-14=", $true, $false, $false, $false, $false, $true, 1, $false, "33+2=", 2) | Out-Null
$find.Execute("99-35=", $true, $false, $false, $false, $false, $true, 1, $false, "70-11=", 2) | Out-Null
$find.Execute("14-0=", $true, $false, $false, $false, $false, $true, 1, $false, "69-26=", 2) | Out-Null
$find.Execute("21+14=", $true, $false, $false, $false, $false, $true, 1, $false, "51+30=", 2) | Out-Null
$find.Execute("73-44=", $true, $false, $false, $false, $false, $true, 1, $false, "44-4=", 2) | Out-Null
$find.Execute("69-53=", $true, $false, $false, $false, $false, $true, 1, $false, "77-29=", 2) | Out-Null
$find.Execute("63+29=", $true, $false, $false, $false, $false, $true, 1, $false, "44+13=", 2) | Out-Null
$find.Execute("33+47=", $true, $false, $false, $false, $false, $true, 1, $false, "29+63=", 2) | Out-Null
$find.Execute("98-27=", $true, $false, $false, $false, $false, $true, 1, $false, "66-22=", 2) | Out-Null
$find.Execute("12-10=", $true, $false, $false, $false, $false, $true, 1, $false, "36+0=", 2) | Out-Null
$find.Execute("52+27=", $true, $false, $false, $false, $false, $true, 1, $false, "80-1=", 2) | Out-Null
$find.Execute("78-38=", $true, $false, $false, $false, $false, $true, 1, $false, "76-59=", 2) | Out-Null
$find.Execute("68+6=", $true, $false, $false, $false, $false, $true, 1, $false, "64+23=", 2) | Out-Null
$find.Execute("45-11=", $true, $false, $false, $false, $false, $true, 1, $false, "32+41=", 2) | Out-Null
$find.Execute("46+22=", $true, $false, $false, $false, $false, $true, 1, $false, "13+66=", 2) | Out-Null
$find.Execute("51-42=", $true, $false, $false, $false, $false, $true, 1, $false, "82-80=", 2) | Out-Null
$find.Execute("11+4=", $true, $false, $false, $false, $false, $true, 1, $false, "47-26=", 2) | Out-Null
$find.Execute("24+74=", $true, $false, $false, $false, $false, $true, 1, $false, "80+13=", 2) | Out-Null
$find.Execute("2+10=", $true, $false, $false, $false, $false, $true, 1, $false, "92-50=", 2) | Out-Null
$find.Execute("30+21=", $true, $false, $false, $false, $false, $true, 1, $false, "15+11=", 2) | Out-Null
$find.Execute("52+33=", $true, $false, $false, $false, $false, $true, 1, $false, "89+6=", 2) | Out-Null
$find.Execute("61+16=", $true, $false, $false, $false, $false, $true, 1, $false, "96-41=", 2) | Out-Null
$find.Execute("61+13=", $true, $false, $false, $false, $false, $true, 1, $false, "21+52=", 2) | Out-Null
$find.Execute("1+76=", $true, $false, $false, $false, $false, $true, 1, $false, "99-91=", 2) | Out-Null
$find.Execute("10+1=", $true, $false, $false, $false, $false, $true, 1, $false, "55-3=", 2) | Out-Null
$find.Execute("38+41=", $true, $false, $false, $false, $false, $true, 1, $false, "14+25=", 2) | Out-Null
$find.Execute("39+3=", $true, $false, $false, $false, $false, $true, 1, $false, "41-21=", 2) | Out-Null
$find.Execute("72+26=", $true, $false, $false, $false, $false, $true, 1, $false, "8+6=", 2) | Out-Null
$find.Execute("1+27=", $true, $false, $false, $false, $false, $true, 1, $false, "52+16=", 2) | Out-Null
$find.Execute("46+16=", $true, $false, $false, $false, $false, $true, 1, $false, "84-69=", 2) | Out-Null
$find.Execute("18-12=", $true, $false, $false, $false, $false, $true, 1, $false, "45+53=", 2) | Out-Null
$find.Execute("7+16=", $true, $false, $false, $false, $false, $true, 1, $false, "61+21=", 2) | Out-Null
$find.Execute("11+81=", $true, $false, $false, $false, $false, $true, 1, $false, "83+7=", 2) | Out-Null
$find.Execute("14+4=", $true, $false, $false, $false, $false, $true, 1, $false, "17+13=", 2) | Out-Null
$find.Execute("57+34=", $true, $false, $false, $false, $false, $true, 1, $false, "30+62=", 2) | Out-Null
$find.Execute("26+43=", $true, $false, $false, $false, $false, $true, 1, $false, "50-18=", 2) | Out-Null
$find.Execute("32-31=", $true, $false, $false, $false, $false, $true, 1, $false, "98-53=", 2) | Out-Null
$find.Execute("72-0=", $true, $false, $false, $false, $false, $true, 1, $false, "26+70=", 2) | Out-Null
$find.Execute("60+20=", $true, $false, $false, $false, $false, $true, 1, $false, "0+67=", 2) | Out-Null
$find.Execute("97-28=", $true, $false, $false, $false, $false, $true, 1, $false, "73-71=", 2) | Out-Null
$find.Execute("29+13=", $true, $false, $false, $false, $false, $true, 1, $false, "41-8=", 2) | Out-Null
$find.Execute("5+72=", $true, $false, $false, $false, $false, $true, 1, $false, "37+44=", 2) | Out-Null
$find.Execute("52-0=", $true, $false, $false, $false, $false, $true, 1, $false, "36-5=", 2) | Out-Null
$find.Execute("4+5=", $true, $false, $false, $false, $false, $true, 1, $false, "83-33=", 2) | Out-Null
$find.Execute("8+59=", $true, $false, $false, $false, $false, $true, 1, $false, "58-14=", 2) | Out-Null
$find.Execute("25+13=", $true, $false, $false, $false, $false, $true, 1, $false, "6+30=", 2) | Out-Null
$find.Execute("87-1=", $true, $false, $false, $false, $false, $true, 1, $false, "28-18=", 2) | Out-Null
$find.Execute("81-17=", $true, $false, $false, $false, $false, $true, 1, $false, "77+3=", 2) | Out-Null
$find.Execute("94-15=", $true, $false, $false, $false, $false, $true, 1, $false, "49+17=", 2) | Out-Null
$find.Execute("42+32=", $true, $false, $false, $false, $false, $true, 1, $false, "7+2=", 2) | Out-Null
$find.Execute("54+17=", $true, $false, $false, $false, $false, $true, 1, $false, "99-74=", 2) | Out-Null
$find.Execute("95-51=", $true, $false, $false, $false, $false, $true, 1, $false, "15+20=", 2) | Out-Null
$find.Execute("74-19=", $true, $false, $false, $false, $false, $true, 1, $false, "46-39=", 2) | Out-Null
$find.Execute("48-4=", $true, $false, $false, $false, $false, $true, 1, $false, "88-70=", 2) | Out-Null
$find.Execute("45+47=", $true, $false, $false, $false, $false, $true, 1, $false, "17+71=", 2) | Out-Null
$find.Execute("89-51=", $true, $false, $false, $false, $false, $true, 1, $false, "48-41=", 2) | Out-Null
$find.Execute("73-37=", $true, $false, $false, $false, $false, $true, 1, $false, "62-41=", 2) | Out-Null
$find.Execute("27+12=", $true, $false, $false, $false, $false, $true, 1, $false, "82-33=", 2) | Out-Null
$find.Execute("4+29=", $true, $false, $false, $false, $false, $true, 1, $false, "56-17=", 2) | Out-Null
$find.Execute("51-2=", $true, $false, $false, $false, $false, $true, 1, $false, "17+47=", 2) | Out-Null
$find.Execute("24+51=", $true, $false, $false, $false, $false, $true, 1, $false, "82-74=", 2) | Out-Null
$find.Execute("9+75=", $true, $false, $false, $false, $false, $true, 1, $false, "63-49=", 2) | Out-Null
$find.Execute("42+19=", $true, $false, $false, $false, $false, $true, 1, $false, "79-43=", 2) | Out-Null
$find.Execute("41+32=", $true, $false, $false, $false, $false, $true, 1, $false, "68-60=", 2) | Out-Null
$find.Execute("36+39=", $true, $false, $false, $false, $false, $true, 1, $false, "30+51=", 2) | Out-Null
$find.Execute("86-11=", $true, $false, $false, $false, $false, $true, 1, $false, "9+56=", 2) | Out-Null
$find.Execute("37+0=", $true, $false, $false, $false, $false, $true, 1, $false, "63+4=", 2) | Out-Null
$find.Execute("64+27=", $true, $false, $false, $false, $false, $true, 1, $false, "75-10=", 2) | Out-Null
$find.Execute("14-4=", $true, $false, $false, $false, $false, $true, 1, $false, "88-2=", 2) | Out-Null
$find.Execute("55+29=", $true, $false, $false, $false, $false, $true, 1, $false, "42-34=", 2) | Out-Null
$find.Execute("41-37=", $true, $false, $false, $false, $false, $true, 1, $false, "29+69=", 2) | Out-Null
$find.Execute("61-14=", $true, $false, $false, $false, $false, $true, 1, $false, "2+45=", 2) | Out-Null
$find.Execute("78-6=", $true, $false, $false, $false, $false, $true, 1, $false, "40+35=", 2) | Out-Null
